# Insert a new data row at row 72 (pushing existing rows 72-164 down to 73-165)
# and populate it with the new "Feria Lagunitas de Puerto Montt - Haba" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(72).Insert()

$ws.Range("A72").Value = 4
$ws.Range("B72").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C72").Value = 'Los Lagos'
$ws.Range("D72").Value = 45195
$ws.Range("E72").Value = 10
$ws.Range("F72").Value = 100112026
$ws.Range("G72").Value = 'Haba'
$ws.Range("H72").Value = 'Sin especificar'
$ws.Range("I72").Value = 'Primera'
$ws.Range("J72").Value = 120
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 18000
$ws.Range("M72").Value = 18000
$ws.Range("N72").Value = '$/saco 25 kilos'
$ws.Range("O72").Value = 'Provincia de Limarí'
$ws.Range("P72").Value = 720
$ws.Range("Q72").Value = 25
$ws.Range("R72").Value = 'Hortaliza'
